$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow columns A:C slightly
$ws.Range("A1:C1").ColumnWidth = 40.42578125

# Add new column M (2023) mirroring column L's formatting
$ws.Range("L3:L12").Copy() | Out-Null
$ws.Range("M3:M12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("M4").Value = 2023
$ws.Range("M5").Value = 311.65582791395695
$ws.Range("M7").Value = 119.55977988994496
$ws.Range("M8").Value = 192.09604802401199
$ws.Range("M10").Value = 78.539269634817401
$ws.Range("M11").Value = 60.030015007503756
$ws.Range("M12").Value = 26.013006503251628

# Row heights
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(5).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 15
$ws.Rows.Item(7).RowHeight = 15
$ws.Rows.Item(8).RowHeight = 15
$ws.Rows.Item(9).RowHeight = 15
$ws.Rows.Item(10).RowHeight = 15
$ws.Rows.Item(11).RowHeight = 15
$ws.Rows.Item(12).RowHeight = 15

# Clear selection on N5, select A1 instead
$ws.Range("A1").Select()
